$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grouped matches")

$data = @(
    ,@('P8211', '{''eft:vidyakaraprabha''}')
    ,@('P8263', '{''eft:leki-d-''}')
    ,@('P8273', '{''eft:rinchen-tso'', ''eft:rin-chen-tsho''}')
    ,@('P8205', '{''eft:yesh-d-ye-shes-sde-'', ''eft:ye-shes-sde'', ''eft:yesh-d-'', ''eft:band-yesh-d-'', ''eft:zhang-yesh-d-'', ''eft:band-yesh-de''}')
    ,@('?', '{''eft:sakyasena''}')
    ,@('P8183', '{''eft:klu-i-rgyal-mtshan'', ''eft:cog-ro-klu-i-rgyal-mtshan''}')
    ,@('P8171', '{''eft:dharmasribhadra''}')
    ,@('P8151', '{''eft:gayadhara''}')
    ,@('P0TMP104', '{''eft:punyasambhava''}')
    ,@('P4259', '{''eft:dpal-gyi-lhun-po'', ''eft:palgyi-lh-npo'', ''eft:ban-de-dpal-gyi-lhun-po''}')
    ,@('P4CZ16780', '{''eft:manjusrigarbha''}')
    ,@('P4258', '{''eft:dpal-byor''}')
    ,@('P3456', '{''eft:tshul-khrims-rgyal-ba''}')
    ,@('P8260', '{''eft:dpal-dbyangs''}')
    ,@('P8213', '{''eft:t-vidyakarasimha'', ''eft:vidyakarasimha''}')
    ,@('P8220', '{''eft:devacandra''}')
    ,@('P0TMP080', '{''eft:hwa-shang-zab-mo''}')
    ,@('P8209', '{''eft:dzi-na-mi-tra-k-'', ''eft:jinamitra'', ''eft:jinamitra-k-''}')
    ,@('P2548', '{''eft:prajnavarma'', ''eft:prajnavarman''}')
    ,@('P5651', '{''eft:pa-tshab-nyi-ma-grags''}')
    ,@('P8268', '{''eft:buddhaprabha''}')
    ,@('P8222', '{''eft:jnanasiddhi'', ''eft:jnanasidhi''}')
    ,@('P4CZ15137', '{''eft:kumarakalasa''}')
    ,@('P8245', '{''eft:buddhakaravarma''}')
    ,@('P1KG8854', '{''eft:srilendrabodhi'', ''eft:silendrabodhi'', ''eft:surendrabodhi''}')
    ,@('P8093', '{''eft:kamalagupta''}')
    ,@('P2637', '{''eft:trakpa-gyaltsen''}')
    ,@('https://lod.dila.edu.tw/resource.php?id=A000089', '{''eft:siladharma''}')
    ,@('P8266', '{''eft:dharmatasila'', ''eft:ch-nyi-tsultrim''}')
    ,@('P8265', '{''eft:ratnaraksita''}')
    ,@('P4CZ16819', '{''eft:sakyaprabha''}')
    ,@('P4263', '{''eft:dge-ba-dpal''}')
    ,@('P0RK8', '{''eft:dharmapala''}')
    ,@('P0TMP092', '{''eft:anandasri-s-''}')
    ,@('P8269', '{''eft:dgon-gling-rma''}')
    ,@('P0TMPT007', '{''eft:rnam-par-mi-rtog-pa''}')
    ,@('P8249', '{''eft:dharmakara''}')
    ,@('P8267', '{''eft:vijayasila''}')
    ,@('P4255', '{''eft:t-jnanagarbha'', ''eft:yesh-nyingpo'', ''eft:ye-shes-snying-po''}')
    ,@('P2956', '{''eft:krsnapandita''}')
    ,@('P8217', '{''eft:t-jnanagarbha'', ''eft:jnanagarbha''}')
    ,@('P3379', '{''eft:dipamkarasrijnana'', ''eft:dipamkara-srijnana''}')
    ,@('P4242', '{''eft:sherab-lekpa''}')
    ,@('P8219', '{''eft:visuddhasimha''}')
    ,@('P8228', '{''eft:surendrabodhi''}')
    ,@('P753', '{''eft:rin-chen-bzang-po''}')
    ,@('P0TMP098', '{''eft:jinavara''}')
    ,@('P3285', '{''eft:sakya-yesh-''}')
    ,@('P8182', '{''eft:dpal-brtsegs'', ''eft:ban-de-dpal-brtsegs'', ''eft:kawa-paltsek-under-the-name-paltsek-raksita-'', ''eft:paltsek'', ''eft:ska-ba-dpal-brtsegs''}')
    ,@('P8206', '{''eft:celu''}')
    ,@('P3709', '{''eft:phakpa-sherab''}')
    ,@('P00KG07267', '{''eft:sarvanyadeva'', ''eft:sarvajnadeva''}')
    ,@('P8261', '{''eft:munivarma'', ''eft:munivarman''}')
    ,@('P3214', '{''eft:danasila''}')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 3).Value = $data[$i][1]
}
